# Auto-generated edit script: rewrite sheet1 A2:I26 to match target state
# - Adds a space before the ":" in the statut_name (col B) entries
# - Fixes two intervention_type (col I) values
# - Reconciles full row contents/order to match the final export
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = '4'
$ws.Cells.Item(2,2).Value = '4 : pas de résultats postés ni publiés'
$ws.Cells.Item(2,3).Value = 'NCT00439582'
$ws.Cells.Item(2,6).Value = '2005'
$ws.Cells.Item(2,7).Value = 'Comparative Effect of 2 Different Sources of Trans Fatty Acid (Milk Fat vs Hydrogenated Oil)on Cardiovascular Risk Factors in Healthy Humans'
$ws.Cells.Item(2,8).Value = 'TRANSFACT1'
$ws.Cells.Item(2,9).Value = 'DIETARY_SUPPLEMENT'

# Row 3
$ws.Cells.Item(3,1).Value = '3'
$ws.Cells.Item(3,2).Value = '3 : résultats postés ou publiés après les 36 mois'
$ws.Cells.Item(3,3).Value = 'NCT00873951'
$ws.Cells.Item(3,6).Value = '2008'
$ws.Cells.Item(3,7).Value = 'Influence of Protein Hydrolysis on Dietary Protein Digestibility and Metabolism in Healthy Subjects'
$ws.Cells.Item(3,9).Value = 'OTHER'

# Row 4
$ws.Cells.Item(4,1).Value = '4'
$ws.Cells.Item(4,2).Value = '4 : pas de résultats postés ni publiés'
$ws.Cells.Item(4,3).Value = 'NCT00685581'
$ws.Cells.Item(4,6).Value = '2008'
$ws.Cells.Item(4,7).Value = 'Rationale, Study Design and Baseline Data of the TRANSQUAL Clinical Trial: A Study to Evaluate the Impact of Different Milk Fatty Acid Profiles on Cardiovascular Risk Factors in Healthy Volunteers; Focus on Trans Fatty Acids'
$ws.Cells.Item(4,8).Value = 'TRANSQUAL WPC'
$ws.Cells.Item(4,9).Value = 'DIETARY_SUPPLEMENT'

# Row 5
$ws.Cells.Item(5,1).Value = '3'
$ws.Cells.Item(5,2).Value = '3 : résultats postés ou publiés après les 36 mois'
$ws.Cells.Item(5,3).Value = 'NCT00617435'
$ws.Cells.Item(5,6).Value = '2009'
$ws.Cells.Item(5,8).Value = 'Trans-Insulin'
$ws.Cells.Item(5,9).Value = 'DIETARY_SUPPLEMENT'

# Row 6
$ws.Cells.Item(6,1).Value = '3'
$ws.Cells.Item(6,2).Value = '3 : résultats postés ou publiés après les 36 mois'
$ws.Cells.Item(6,3).Value = 'NCT00931151'
$ws.Cells.Item(6,6).Value = '2009'
$ws.Cells.Item(6,9).Value = 'DIETARY_SUPPLEMENT'

# Row 7
$ws.Cells.Item(7,1).Value = '4'
$ws.Cells.Item(7,2).Value = '4 : pas de résultats postés ni publiés'
$ws.Cells.Item(7,3).Value = 'NCT00690781'
$ws.Cells.Item(7,6).Value = '2010'
$ws.Cells.Item(7,7).Value = 'Effect of Milk Proteins and Protein Feeding Pattern on Body Composition and Protein Metabolism in Energy Restricted Obese Subjects'
$ws.Cells.Item(7,8).Value = 'SURPROL-CF-H'
$ws.Cells.Item(7,9).Value = 'DIETARY_SUPPLEMENT'

# Row 8
$ws.Cells.Item(8,1).Value = '4'
$ws.Cells.Item(8,2).Value = '4 : pas de résultats postés ni publiés'
$ws.Cells.Item(8,3).Value = 'NCT01209572'
$ws.Cells.Item(8,6).Value = '2010'
$ws.Cells.Item(8,7).Value = 'Modelling of 24h Energy Expenditure From Heart Rate, Actimetry and Other Parameters Recorded Under Free-living Conditions'
$ws.Cells.Item(8,8).Value = 'Modelheart'
$ws.Cells.Item(8,9).Value = 'DEVICE'

# Row 9
$ws.Cells.Item(9,1).Value = '3'
$ws.Cells.Item(9,2).Value = '3 : résultats postés ou publiés après les 36 mois'
$ws.Cells.Item(9,3).Value = 'NCT00994526'
$ws.Cells.Item(9,6).Value = '2010'
$ws.Cells.Item(9,7).Value = 'Effect of Processed Meat on Colorectal Carcinogenesis. Study of Mechanisms. Choice of Preventive Strategies'
$ws.Cells.Item(9,8).Value = 'Hemcancer'
$ws.Cells.Item(9,9).Value = 'OTHER'

# Row 10
$ws.Cells.Item(10,1).Value = '2'
$ws.Cells.Item(10,2).Value = '2 : résultats postés ou publiés entre 12 et 36 mois'
$ws.Cells.Item(10,3).Value = 'NCT00862329'
$ws.Cells.Item(10,6).Value = '2010'
$ws.Cells.Item(10,7).Value = ''
$ws.Cells.Item(10,8).Value = ''
$ws.Cells.Item(10,9).Value = 'DIETARY_SUPPLEMENT'

# Row 11
$ws.Cells.Item(11,1).Value = '4'
$ws.Cells.Item(11,2).Value = '4 : pas de résultats postés ni publiés'
$ws.Cells.Item(11,3).Value = 'NCT01995162'
$ws.Cells.Item(11,6).Value = '2013'
$ws.Cells.Item(11,7).Value = 'A Smartphone Application to Evaluate Energy Expenditure and Duration of Moderate-intensity Activities in Free-living Conditions (eMouve 2)'
$ws.Cells.Item(11,9).Value = 'DEVICE'

# Row 12
$ws.Cells.Item(12,1).Value = '4'
$ws.Cells.Item(12,2).Value = '4 : pas de résultats postés ni publiés'
$ws.Cells.Item(12,3).Value = 'NCT01995253'
$ws.Cells.Item(12,6).Value = '2013'
$ws.Cells.Item(12,7).Value = 'A Smartphone Application to Evaluate Energy Expenditure and Duration of Moderate-intensity Activities'
$ws.Cells.Item(12,9).Value = 'DEVICE'

# Row 13
$ws.Cells.Item(13,1).Value = '4'
$ws.Cells.Item(13,2).Value = '4 : pas de résultats postés ni publiés'
$ws.Cells.Item(13,3).Value = 'NCT02348554'
$ws.Cells.Item(13,6).Value = '2014'
$ws.Cells.Item(13,7).Value = 'A Smartphone Application to Evaluate Energy Expenditure and Duration of Activities in Free-living Conditions for Overweight and Obese People (eMouve3)'
$ws.Cells.Item(13,8).Value = 'eMouve3'
$ws.Cells.Item(13,9).Value = 'BEHAVIORAL'

# Row 14
$ws.Cells.Item(14,1).Value = '2'
$ws.Cells.Item(14,2).Value = '2 : résultats postés ou publiés entre 12 et 36 mois'
$ws.Cells.Item(14,3).Value = 'NCT02354794'
$ws.Cells.Item(14,6).Value = '2014'
$ws.Cells.Item(14,7).Value = 'Effect of Oral Supplementation With One Form of L-arginine on Vascular Endothelial Function in Healthy Subjects Featuring Risk Factors Related to the Metabolic Syndrome.'
$ws.Cells.Item(14,9).Value = 'DIETARY_SUPPLEMENT'

# Row 15
$ws.Cells.Item(15,1).Value = '3'
$ws.Cells.Item(15,2).Value = '3 : résultats postés ou publiés après les 36 mois'
$ws.Cells.Item(15,3).Value = 'NCT02157805'
$ws.Cells.Item(15,6).Value = '2014'
$ws.Cells.Item(15,7).Value = 'Effect of Technological Processes on Nutritional Quality of Meat Proteins'
$ws.Cells.Item(15,9).Value = 'OTHER'

# Row 16
$ws.Cells.Item(16,1).Value = '4'
$ws.Cells.Item(16,2).Value = '4 : pas de résultats postés ni publiés'
$ws.Cells.Item(16,3).Value = 'NCT02473302'
$ws.Cells.Item(16,6).Value = '2014'
$ws.Cells.Item(16,7).Value = 'Preventive Strategies in Colorectal Carcinogenesis Production and Meat Processing'
$ws.Cells.Item(16,8).Value = ''
$ws.Cells.Item(16,9).Value = 'OTHER'

# Row 17
$ws.Cells.Item(17,1).Value = '2'
$ws.Cells.Item(17,2).Value = '2 : résultats postés ou publiés entre 12 et 36 mois'
$ws.Cells.Item(17,3).Value = 'NCT03492593'
$ws.Cells.Item(17,6).Value = '2016'
$ws.Cells.Item(17,7).Value = 'Métabolismes Des caroténoïdes Dans la lumière du Tube Digestif de l''Homme Sain'
$ws.Cells.Item(17,8).Value = 'CarotenoiDig'
$ws.Cells.Item(17,9).Value = 'OTHER'

# Row 18
$ws.Cells.Item(18,1).Value = '3'
$ws.Cells.Item(18,2).Value = '3 : résultats postés ou publiés après les 36 mois'
$ws.Cells.Item(18,3).Value = 'NCT03265392'
$ws.Cells.Item(18,6).Value = '2018'
$ws.Cells.Item(18,7).Value = 'Digestion: Building a Better Health and Better Understanding the Digestive Processes Thanks to Magnetic Resonance Imaging'
$ws.Cells.Item(18,8).Value = 'DECOUVRIR-M'
$ws.Cells.Item(18,9).Value = 'OTHER'

# Row 19
$ws.Cells.Item(19,1).Value = '2'
$ws.Cells.Item(19,2).Value = '2 : résultats postés ou publiés entre 12 et 36 mois'
$ws.Cells.Item(19,3).Value = 'NCT03279211'
$ws.Cells.Item(19,6).Value = '2019'
$ws.Cells.Item(19,7).Value = 'True Ileal Amino Acid Digestibility of Whey and Zein Proteins in Healthy Volunteers With Naso-ileal Tubes'
$ws.Cells.Item(19,9).Value = 'OTHER'

# Row 20
$ws.Cells.Item(20,1).Value = '2'
$ws.Cells.Item(20,2).Value = '2 : résultats postés ou publiés entre 12 et 36 mois'
$ws.Cells.Item(20,3).Value = 'NCT04072770'
$ws.Cells.Item(20,6).Value = '2020'
$ws.Cells.Item(20,7).Value = 'Bioavailability of Protein and Amino Acids of Pea Protein Isolate in Healthy Volunteers'
$ws.Cells.Item(20,8).Value = 'Qualipois'
$ws.Cells.Item(20,9).Value = 'DIETARY_SUPPLEMENT'

# Row 21
$ws.Cells.Item(21,1).Value = '4'
$ws.Cells.Item(21,2).Value = '4 : pas de résultats postés ni publiés'
$ws.Cells.Item(21,3).Value = 'NCT06624033'
$ws.Cells.Item(21,6).Value = '2023'
$ws.Cells.Item(21,7).Value = 'Single-blind, Randomized, Cross-over Comparative Bioavailability Study About the Kinetics of Plasma Amino Acid Concentrations Subsequent to the Consumption of Innovative Legume-based Products.'
$ws.Cells.Item(21,8).Value = 'LEG''UP'
$ws.Cells.Item(21,9).Value = 'OTHER'

# Row 22
$ws.Cells.Item(22,1).Value = '4'
$ws.Cells.Item(22,2).Value = '4 : pas de résultats postés ni publiés'
$ws.Cells.Item(22,3).Value = 'NCT05047757'
$ws.Cells.Item(22,6).Value = '2023'
$ws.Cells.Item(22,7).Value = 'Fava Bean Protein and Amino Acid Bioavailability in Healthy Volunteers'
$ws.Cells.Item(22,8).Value = 'Leg4Life'
$ws.Cells.Item(22,9).Value = 'DIETARY_SUPPLEMENT'

# Row 23
$ws.Cells.Item(23,1).Value = '3'
$ws.Cells.Item(23,2).Value = '3 : résultats postés ou publiés après les 36 mois'
$ws.Cells.Item(23,3).Value = 'NCT00862017'
$ws.Cells.Item(23,7).Value = 'Effect of Monosodium Glutamate on Gastric Emptying and Postprandial Nitrogen in Healthy Volunteers'
$ws.Cells.Item(23,9).Value = 'DIETARY_SUPPLEMENT'

# Row 24
$ws.Cells.Item(24,1).Value = '4'
$ws.Cells.Item(24,2).Value = '4 : pas de résultats postés ni publiés'
$ws.Cells.Item(24,3).Value = 'NCT01154608'
$ws.Cells.Item(24,7).Value = ''
$ws.Cells.Item(24,9).Value = 'OTHER'

# Row 25
$ws.Cells.Item(25,1).Value = '3'
$ws.Cells.Item(25,2).Value = '3 : résultats postés ou publiés après les 36 mois'
$ws.Cells.Item(25,3).Value = 'NCT01154582'
$ws.Cells.Item(25,7).Value = ''
$ws.Cells.Item(25,9).Value = 'DIETARY_SUPPLEMENT'

# Row 26
$ws.Cells.Item(26,1).Value = '3'
$ws.Cells.Item(26,2).Value = '3 : résultats postés ou publiés après les 36 mois'
$ws.Cells.Item(26,3).Value = 'NCT02352740'
$ws.Cells.Item(26,7).Value = 'Characterization of the Metabolic Fate of an Oral L-arginine Form in Healthy Subjects Featuring Risk Factors Related to the Metabolic Syndrome.'
$ws.Cells.Item(26,9).Value = 'DIETARY_SUPPLEMENT'

